$wb = $excel.ActiveWorkbook

# The "Indennità" column header lives on the "modello_del_foglio_di_dati" sheet,
# in cell L1. We shorten the descriptive sub-text by removing the trailing
# ", 1/12 del totale annuo" clause (that clause is now only used by the
# "Pagamenti straordinari" header elsewhere in the sheet).
$ws = $wb.Worksheets.Item("modello_del_foglio_di_dati")

$cell = $ws.Range("L1")
$boldPart = "Indennità"
$restPart = [char]10 + "(lavoro a turni, lavoro domenicale e notturno nonché altre indennità per faticosità del lavoro)"
$text = $boldPart + $restPart
$cell.Value = $text

# Preserve the same rich-text run layout used by the sibling headers in this
# sheet: a bold title run followed by a regular-weight description run.
$titleRun = $cell.Characters(1, $boldPart.Length)
$titleRun.Font.Bold = $true
$titleRun.Font.Size = 9
$titleRun.Font.Name = "Arial"

$descRun = $cell.Characters($boldPart.Length + 1, $restPart.Length)
$descRun.Font.Bold = $false
$descRun.Font.Size = 9
$descRun.Font.Name = "Arial"
$descRun.Font.Color = 0
